$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated nombre_aides (column C) and montant_total (column D) values
# for rows reflecting the 2020-06-25 data refresh.
$updates = @(
    @(2, 277715, 355121162),
    @(8, 769, 1135394),
    @(10, 105988, 155584022),
    @(12, 52214, 75473447),
    @(16, 3429, 4874603),
    @(20, 4908, 6874088),
    @(22, 68858, 86478690),
    @(26, 10, 15000),
    @(27, 257, 370718),
    @(28, 29640, 43449773),
    @(30, 10179, 14687443),
    @(33, 1351, 1900715),
    @(35, 1404, 1979703),
    @(36, 87144, 110495459),
    @(37, 57, 65593),
    @(38, 74, 104153),
    @(42, 804, 1183958),
    @(43, 8, 12000),
    @(44, 40775, 59878093),
    @(46, 8155, 11713027),
    @(48, 1186, 1647017),
    @(51, 1729, 2411827),
    @(52, 60536, 76225032),
    @(56, 340, 500909),
    @(58, 25644, 37651295),
    @(61, 9837, 14243937),
    @(62, 1198, 1674607),
    @(65, 1115, 1566576),
    @(67, 17776, 23270436),
    @(71, 6465, 9459992),
    @(72, 4407, 6405710),
    @(73, 405, 568758),
    @(74, 219, 308056),
    @(75, 123920, 155152951),
    @(78, 18, 27000),
    @(79, 366, 534782),
    @(81, 57653, 84645609),
    @(82, 69, 102082),
    @(84, 26414, 38250161),
    @(86, 2305, 3324328),
    @(87, 2100, 2956207),
    @(88, 23784, 32297817),
    @(92, 6356, 9368176),
    @(94, 5523, 8006034),
    @(96, 386, 545809),
    @(97, 321, 464161),
    @(98, 5394, 7473286),
    @(100, 1368, 2014062),
    @(102, 1871, 2730601),
    @(104, 65, 93020),
    @(105, 95, 133789),
    @(106, 125128, 155153237),
    @(110, 864, 1270632),
    @(112, 48250, 70826657),
    @(113, 79, 117659),
    @(114, 23819, 34533935),
    @(115, 1102, 1502832),
    @(118, 1709, 2404872),
    @(120, 37421, 50207259),
    @(121, 23, 32653),
    @(126, 12429, 18273127),
    @(127, 3247, 4688941),
    @(130, 326, 472036),
    @(132, 295, 413190),
    @(133, 14150, 18777683),
    @(137, 6071, 8848752),
    @(139, 4020, 5799615),
    @(141, 212, 289433),
    @(142, 173, 250664),
    @(144, 8872, 12805026),
    @(145, 1309, 1945078),
    @(146, 180, 265671),
    @(147, 26, 38690),
    @(149, 77620, 97513434),
    @(154, 589, 870037),
    @(156, 31064, 45606600),
    @(158, 11501, 16637592),
    @(160, 1055, 1482396),
    @(162, 1253, 1774226),
    @(164, 209991, 261928662),
    @(166, 151, 218453),
    @(170, 788, 1158790),
    @(172, 79169, 116157136),
    @(175, 29377, 42316015),
    @(178, 4418, 6310386),
    @(181, 3722, 5178659),
    @(184, 229021, 284180113),
    @(185, 138, 151352),
    @(186, 219, 311911),
    @(190, 556, 809382),
    @(192, 85888, 125786042),
    @(195, 45307, 65569267),
    @(198, 4019, 5654667),
    @(201, 4320, 5981557),
    @(204, 94085, 118431412),
    @(205, 68, 70104),
    @(209, 508, 742412),
    @(211, 45479, 66697294),
    @(213, 10945, 15750706),
    @(215, 1687, 2417347),
    @(216, 8, 10096),
    @(217, 1948, 2719580),
    @(218, 225212, 285059174),
    @(219, 149, 185356),
    @(220, 225, 324832),
    @(224, 734, 1080810),
    @(226, 86433, 126808751),
    @(229, 57138, 82914285),
    @(231, 2082, 2942645),
    @(234, 3542, 4972531),
    @(235, 387829, 509494954),
    @(236, 70, 90997),
    @(240, 1152, 1709551),
    @(242, 179221, 263875195),
    @(243, 318, 473790),
    @(245, 149596, 217607407),
    @(247, 2239, 3141648),
    @(249, 4586, 6471254),
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 3).Value = $u[1]
    $ws.Cells.Item($r, 4).Value = $u[2]
}

Write-Output "Applied $($updates.Count) row updates"
